$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "a" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "a"

# Header row 1 (repeated "EU27+UK" label across C1:F1) - first use of this
# string anchors its shared-string index.
$ws.Range("C1").Value = "EU27+UK"
$ws.Range("D1").Value = "EU27+UK"
$ws.Range("E1").Value = "EU27+UK"
$ws.Range("F1").Value = "EU27+UK"

# Row labels - set the four region names in the order they were first
# typed by the author (EU27+UK, China, USA, RoW) so shared-string indices
# line up, before filling in the rest of the table.
$ws.Range("A3").Value = "China"
$ws.Range("A15").Value = "USA"
$ws.Range("A11").Value = "RoW"

# Header row 2 - process labels, centered like the existing tables
$ws.Range("C2").Value = "Refinery of Generators of Onshore Wind Turbines"
$ws.Range("D2").Value = "Refinery of Generators of Offshore Wind Turbines"
$ws.Range("E2").Value = "Refinery of Silicon layer in PV panel"
$ws.Range("F2").Value = "Refinery of Cu in wires of WT and PV"
$ws.Range("C2:F2").HorizontalAlignment = -4108

# Material labels in column B, centered
$materials = "Neodymium", "Dysprosium", "Copper ores and concentrates", "Raw silicon"
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item(3 + $i, 2).Value = $materials[$i]
    $ws.Cells.Item(7 + $i, 2).Value = $materials[$i]
    $ws.Cells.Item(11 + $i, 2).Value = $materials[$i]
    $ws.Cells.Item(15 + $i, 2).Value = $materials[$i]
}
$ws.Range("B3:B18").HorizontalAlignment = -4108

# Region labels filling out column A for every block of 4 rows
$ws.Range("A3").Value = "China"
$ws.Range("A4").Value = "China"
$ws.Range("A5").Value = "China"
$ws.Range("A6").Value = "China"
$ws.Range("A7").Value = "EU27+UK"
$ws.Range("A8").Value = "EU27+UK"
$ws.Range("A9").Value = "EU27+UK"
$ws.Range("A10").Value = "EU27+UK"
$ws.Range("A11").Value = "RoW"
$ws.Range("A12").Value = "RoW"
$ws.Range("A13").Value = "RoW"
$ws.Range("A14").Value = "RoW"
$ws.Range("A15").Value = "USA"
$ws.Range("A16").Value = "USA"
$ws.Range("A17").Value = "USA"
$ws.Range("A18").Value = "USA"

# Allocation matrix values (C:F, rows 3-18)
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.792
$ws.Range("F6").Value = 0

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0

$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0

$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1

$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.08
$ws.Range("F10").Value = 0

$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0

$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0

$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0

$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.072
$ws.Range("F14").Value = 0

$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0

$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0

$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0

$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.056
$ws.Range("F18").Value = 0

# Column widths (approximate the author's bestFit auto-sizing)
$ws.Columns("B").ColumnWidth = 25.333333333333336
$ws.Columns("C").ColumnWidth = 42.0
$ws.Columns("D").ColumnWidth = 42.33333333333333
$ws.Columns("E").ColumnWidth = 29.5
$ws.Columns("F").ColumnWidth = 30.666666666666664

# Selection left by the author on this sheet
[void]$ws.Range("D33").Select()

# --- Adjust the "S" sheet (selection + column A width) ---
$sSheet = $wb.Worksheets.Item("S")
$sSheet.Columns("A").ColumnWidth = 42.33333333333333
[void]$sSheet.Range("A4:A7").Select()

# --- Adjust the "RR" sheet selection (it was the previously active tab) ---
$rrSheet = $wb.Worksheets.Item("RR")
[void]$rrSheet.Range("B1:E1").Select()

# --- Make the new sheet the active tab, matching the author's final view ---
[void]$ws.Select()
